# User guide / sample data update for co2mpas_driver model
$wb = $excel.ActiveWorkbook

# --- "inputs" sheet: desired_velocity 10 -> 40 ---
$wsInputs = $wb.Worksheets.Item("inputs")
$wsInputs.Range("B2").Value = 40

# --- "config" sheet: vehicle_id 35135 -> 39393 ---
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Range("B2").Value = 39393

# --- "vehicle_inputs" sheet: no data changes ---
$wsVehicle = $wb.Worksheets.Item("vehicle_inputs")

# --- "time_series" sheet: add a "velocities" column (B) ---
$wsSeries = $wb.Worksheets.Item("time_series")
$wsSeries.Range("B2").Value = 5
$wsSeries.Range("B3").Formula = "=B2 + 5"
$wsSeries.Range("B4:B22").Formula = "=B3 + 5"

# --- restore per-sheet selections (cosmetic UI state) ---
$wsInputs.Range("D24").Select()
$wsConfig.Range("C11").Select()
$wsVehicle.Range("F20").Select()
$wsSeries.Range("F15").Select()
